{"js": "// Indent the \"Output:\" block of the \"Division of zero\" scenario (Question 3)\n// by adding a 0.5\" (36pt / 720 twips) left indent to each paragraph of\n// console-output text that follows that particular \"Output:\" heading.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Division of zero\" scenario heading.\nlet scenarioIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"Division of zero\") {\n    scenarioIdx = i;\n    break;\n  }\n}\nif (scenarioIdx === -1) {\n  throw new Error('Could not find the \"Division of zero\" scenario heading.');\n}\n\n// From the scenario heading, find the next \"Output:\" label paragraph.\nlet outputIdx = -1;\nfor (let i = scenarioIdx + 1; i < items.length; i++) {\n  if (items[i].text.trim() === \"Output:\") {\n    outputIdx = i;\n    break;\n  }\n}\nif (outputIdx === -1) {\n  throw new Error('Could not find the \"Output:\" paragraph for the scenario.');\n}\n\n// The output block runs from right after \"Output:\" until the first blank\n// paragraph (which separates this scenario from the next one).\nconst targets = [];\nfor (let i = outputIdx + 1; i < items.length; i++) {\n  const text = items[i].text.trim();\n  if (text === \"\") break;\n  targets.push(items[i]);\n}\n\n// Apply a 36pt (720 twips / 0.5in) left indent to each output-block paragraph.\nfor (const p of targets) {\n  p.leftIndent = 36;\n}\n\nawait context.sync();\n", "ps1": "# Indent the \"Output:\" block of the \"Division of zero\" scenario (Question 3)\n# by adding a 0.5\" (36pt / 720 twips) left indent to each paragraph of\n# console-output text that follows that particular \"Output:\" heading.\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n# Locate the \"Division of zero\" scenario heading.\n$scenarioIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"Division of zero\") {\n        $scenarioIdx = $i\n        break\n    }\n}\nif ($scenarioIdx -eq -1) {\n    throw \"Could not find the 'Division of zero' scenario heading.\"\n}\n\n# From the scenario heading, find the next \"Output:\" label paragraph.\n$outputIdx = -1\nfor ($i = $scenarioIdx + 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"Output:\") {\n        $outputIdx = $i\n        break\n    }\n}\nif ($outputIdx -eq -1) {\n    throw \"Could not find the 'Output:' paragraph for the scenario.\"\n}\n\n# The output block runs from right after \"Output:\" until the first blank\n# paragraph (which separates this scenario from the next one). Indent each\n# paragraph in that block by 36pt (720 twips / 0.5in).\nfor ($i = $outputIdx + 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq \"\") { break }\n    $p.Range.ParagraphFormat.LeftIndent = 36\n}\n"}
